# Insert a new data row at row 36 (shifts existing rows 36-64 down to 37-65)
# and populate it with the new record reported for this market/product.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:36").Insert()

$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value = "Bíobío"
$ws.Cells.Item(36, 4).Value = 44790
$ws.Cells.Item(36, 5).Value = 8
$ws.Cells.Item(36, 6).Value = 100112037
$ws.Cells.Item(36, 7).Value = "Cebollín"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 300
$ws.Cells.Item(36, 11).Value = 5000
$ws.Cells.Item(36, 12).Value = 6000
$ws.Cells.Item(36, 13).Value = 5333
$ws.Cells.Item(36, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 148
$ws.Cells.Item(36, 17).Value = 36
$ws.Cells.Item(36, 18).Value = "Hortaliza"
